$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 21153.4
$ws.Range("I47").Value = 5767
$ws.Range("J47").Value = 25000
$ws.Range("K47").Value = 5767
$ws.Range("L47").Value = 25000
$ws.Range("M47").Value = -4795
$ws.Range("N47").Value = -26944
$ws.Range("H54").Value = 10010
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H101").Value = 7938750
$ws.Range("I101").Value = 10991150
$ws.Range("J101").Value = 2510.6
$ws.Range("K101").Value = 32973450
$ws.Range("L101").Value = 7531.799999999999
$ws.Range("M101").Value = -32971828
$ws.Range("N101").Value = -10775.8
$ws.Range("H113").Value = 13998.637
$ws.Range("I113").Value = 16499.166
$ws.Range("J113").Value = 10998
$ws.Range("K113").Value = 16499.166
$ws.Range("L113").Value = 10998
$ws.Range("M113").Value = -13245.166
$ws.Range("N113").Value = -17506
$ws.Range("H125").Value = 8318.286
$ws.Range("I125").Value = 17666
$ws.Range("K125").Value = 158994
$ws.Range("M125").Value = -156534
$ws.Range("H129").Value = 1397.25
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H137").Value = 7607.853
$ws.Range("I137").Value = 9234.879999999999
$ws.Range("J137").Value = 3088.3333
$ws.Range("K137").Value = 27704.64
$ws.Range("L137").Value = 9264.999899999999
$ws.Range("M137").Value = -25154.64
$ws.Range("N137").Value = -14364.9999
$ws.Range("H138").Value = 2247.6843
$ws.Range("I138").Value = 900.75
$ws.Range("J138").Value = 3459.925
$ws.Range("K138").Value = 2702.25
$ws.Range("L138").Value = 10379.775
$ws.Range("M138").Value = 2437.75
$ws.Range("N138").Value = -20659.775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2872063.8
$ws.Range("J6").Value = 25000.5
$ws.Range("L6").Value = 25000.5
$ws.Range("N6").Value = -25346.5
$ws.Range("H61").Value = 4187.357
$ws.Range("I61").Value = 3880.652
$ws.Range("K61").Value = 3880.652
$ws.Range("M61").Value = -3668.652
$ws.Range("H74").Value = 2485.0476
$ws.Range("I74").Value = 1508.3077
$ws.Range("K74").Value = 1508.3077
$ws.Range("M74").Value = -634.3077000000001
$ws.Range("H77").Value = 2485.0476
$ws.Range("I77").Value = 1508.3077
$ws.Range("K77").Value = 7541.538500000001
$ws.Range("M77").Value = -3173.538500000001
$ws.Range("H102").Value = 8875.395
$ws.Range("I102").Value = 10459.521
$ws.Range("K102").Value = 10459.521
$ws.Range("M102").Value = -8837.521000000001
$ws.Range("H105").Value = 370000
$ws.Range("J105").Value = 370000
$ws.Range("L105").Value = 370000
$ws.Range("N105").Value = -376988
$ws.Range("H136").Value = 4187.357
$ws.Range("I136").Value = 3880.652
$ws.Range("K136").Value = 11641.956
$ws.Range("M136").Value = -9091.956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 8641709
$ws.Range("I7").Value = 6415.8335
$ws.Range("J7").Value = 19004060
$ws.Range("K7").Value = 6415.8335
$ws.Range("L7").Value = 19004060
$ws.Range("M7").Value = -6302.8335
$ws.Range("N7").Value = -19004286
$ws.Range("H20").Value = 3072
$ws.Range("I20").Value = 1882.7858
$ws.Range("K20").Value = 1882.7858
$ws.Range("M20").Value = -1635.7858
$ws.Range("H22").Value = 39
$ws.Range("I22").Value = 39
$ws.Range("K22").Value = 39
$ws.Range("M22").Value = 134
$ws.Range("H107").Value = 2573.1365
$ws.Range("J107").Value = 3248.5
$ws.Range("L107").Value = 3248.5
$ws.Range("N107").Value = -7088.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 500
$ws.Range("J12").Value = 500
$ws.Range("L12").Value = 500
$ws.Range("N12").Value = -840
$ws.Range("H22").Value = 1100.8077
$ws.Range("I22").Value = 707.53845
$ws.Range("K22").Value = 707.53845
$ws.Range("M22").Value = -357.53845
$ws.Range("H31").Value = 7821.1113
$ws.Range("I31").Value = 8341.522999999999
$ws.Range("J31").Value = 5999.6665
$ws.Range("K31").Value = 8341.522999999999
$ws.Range("L31").Value = 5999.6665
$ws.Range("M31").Value = -8046.522999999999
$ws.Range("N31").Value = -6589.6665
$ws.Range("H34").Value = 7821.1113
$ws.Range("I34").Value = 8341.522999999999
$ws.Range("J34").Value = 5999.6665
$ws.Range("K34").Value = 8341.522999999999
$ws.Range("L34").Value = 5999.6665
$ws.Range("M34").Value = -8139.522999999999
$ws.Range("N34").Value = -6403.6665
$ws.Range("H58").Value = 3010.6956
$ws.Range("I58").Value = 2842.35
$ws.Range("J58").Value = 4133
$ws.Range("K58").Value = 2842.35
$ws.Range("L58").Value = 4133
$ws.Range("M58").Value = -2639.35
$ws.Range("N58").Value = -4539
$ws.Range("H130").Value = 54950
$ws.Range("J130").Value = 54950
$ws.Range("L130").Value = 54950
$ws.Range("N130").Value = -64990
$ws.Range("H134").Value = 2703.611
$ws.Range("I134").Value = 3143.0833
$ws.Range("J134").Value = 1824.6666
$ws.Range("K134").Value = 9429.249899999999
$ws.Range("L134").Value = 5473.9998
$ws.Range("M134").Value = -6894.249899999999
$ws.Range("N134").Value = -10543.9998
$ws.Range("H136").Value = 3010.6956
$ws.Range("I136").Value = 2842.35
$ws.Range("J136").Value = 4133
$ws.Range("K136").Value = 8527.049999999999
$ws.Range("L136").Value = 12399
$ws.Range("M136").Value = -5977.049999999999
$ws.Range("N136").Value = -17499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 123.4
$ws.Range("J12").Value = 29.4
$ws.Range("L12").Value = 88.19999999999999
$ws.Range("N12").Value = -434.2
$ws.Range("H114").Value = 3282
$ws.Range("I114").Value = 1095.3334
$ws.Range("J114").Value = 4594
$ws.Range("K114").Value = 3286.0002
$ws.Range("L114").Value = 13782
$ws.Range("M114").Value = -32.00019999999995
$ws.Range("N114").Value = -20290
$ws.Range("H117").Value = 2188.077
$ws.Range("I117").Value = 3205.25
$ws.Range("J117").Value = 1736
$ws.Range("K117").Value = 9615.75
$ws.Range("L117").Value = 5208
$ws.Range("M117").Value = -6173.75
$ws.Range("N117").Value = -12092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4279747.5
$ws.Range("J70").Value = 8853.857
$ws.Range("L70").Value = 8853.857
$ws.Range("N70").Value = -9393.857
$ws.Range("H73").Value = 4279747.5
$ws.Range("J73").Value = 8853.857
$ws.Range("L73").Value = 8853.857
$ws.Range("N73").Value = -10725.857
$ws.Range("H80").Value = 7398.2856
$ws.Range("I80").Value = 7398.2856
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 7398.2856
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -6400.2856
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 7398.2856
$ws.Range("I83").Value = 7398.2856
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 36991.428
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -31999.428
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 8106.6113
$ws.Range("I97").Value = 10066.857
$ws.Range("K97").Value = 10066.857
$ws.Range("M97").Value = -9570.857
$ws.Range("H102").Value = 9569
$ws.Range("I102").Value = 10660.571
$ws.Range("K102").Value = 10660.571
$ws.Range("M102").Value = -9038.571
$ws.Range("H106").Value = 455750
$ws.Range("J106").Value = 455750
$ws.Range("L106").Value = 455750
$ws.Range("N106").Value = -458274
$ws.Range("H132").Value = 2651.4614
$ws.Range("I132").Value = 2497.4167
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 7492.250100000001
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -4962.250100000001
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7619.7617
$ws.Range("I16").Value = 8702.5
$ws.Range("J16").Value = 4155
$ws.Range("K16").Value = 8702.5
$ws.Range("L16").Value = 4155
$ws.Range("M16").Value = -8532.5
$ws.Range("N16").Value = -4495
$ws.Range("H93").Value = 6418.636
$ws.Range("I93").Value = 7042.1763
$ws.Range("J93").Value = 4298.6
$ws.Range("K93").Value = 7042.1763
$ws.Range("L93").Value = 4298.6
$ws.Range("M93").Value = -5794.1763
$ws.Range("N93").Value = -6794.6
$ws.Range("H105").Value = 45750
$ws.Range("J105").Value = 45750
$ws.Range("L105").Value = 45750
$ws.Range("N105").Value = -52738
$ws.Range("H136").Value = 4711.077
$ws.Range("I136").Value = 2795.2856
$ws.Range("J136").Value = 6946.1665
$ws.Range("K136").Value = 8385.856800000001
$ws.Range("L136").Value = 20838.4995
$ws.Range("M136").Value = -5835.856800000001
$ws.Range("N136").Value = -25938.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 20550.182
$ws.Range("J126").Value = 5257.2
$ws.Range("L126").Value = 15771.6
$ws.Range("N126").Value = -20711.6
$ws.Range("H132").Value = 10413.761
$ws.Range("I132").Value = 11000.878
$ws.Range("J132").Value = 5599.4
$ws.Range("K132").Value = 33002.63400000001
$ws.Range("L132").Value = 16798.2
$ws.Range("M132").Value = -30472.63400000001
$ws.Range("N132").Value = -21858.2
$ws.Range("H136").Value = 1407076
$ws.Range("I136").Value = 2573917.2
$ws.Range("J136").Value = 6866.4
$ws.Range("K136").Value = 7721751.600000001
$ws.Range("L136").Value = 20599.2
$ws.Range("M136").Value = -7719201.600000001
$ws.Range("N136").Value = -25699.2
